$d = $word.ActiveDocument

# 1. Simplify the title paragraph: merge the three runs ("Melt flow index (MFI),
#    Shore hardness, " + "Vicat" + "/HDT and LOI analysis") into a single run and
#    drop the spell-check proofErr markers. Find/Replace with identical
#    before/after text (but spanning multiple runs) causes Word to consolidate
#    the matched text into one run using the first run's formatting and to clear
#    the now-stale proofErr tags.
$d.Content.Find.Execute(
    "Melt flow index (MFI), Shore hardness, Vicat/HDT and LOI analysis",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Melt flow index (MFI), Shore hardness, Vicat/HDT and LOI analysis", 2) | Out-Null

# 2. Remove the stray _GoBack bookmark that currently sits at the end of the
#    title paragraph -- it needs to move down to the new name paragraph below.
$d.Bookmarks.Item("_GoBack").Delete()

# 3. Add a new centered name line "Rosa Ilaria Quercia 207014" right after
#    "Gabriel Orsi 207696", matching the formatting of the surrounding name
#    paragraphs (Helvetica, bCs, 14pt/28 half-points).
$gabrielPara = $d.Paragraphs.Item(19)
$gabrielPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(20)
# A trailing placeholder character is appended and removed below; it sidesteps
# a quirk where adding a bookmark collapsed exactly on the last character
# before a paragraph mark gets mis-anchored.
$newPara.Range.Text = "Rosa Ilaria Quercia 207014X"

# 4. Re-create the _GoBack bookmark right after the new text (collapsed, same
#    placement pattern it originally had after the title).
$newPara2 = $d.Paragraphs.Item(20)
$endWithPlaceholder = $newPara2.Range.End
$bmPos = $endWithPlaceholder - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 5. Drop the placeholder character now that the bookmark is safely anchored.
$placeholderRange = $d.Range($endWithPlaceholder - 2, $endWithPlaceholder - 1)
$placeholderRange.Delete()
